$wb = $excel.ActiveWorkbook

# The "HOC" worksheet holds the single opportunity-cost figure this edit updates.
$ws = $wb.Worksheets.Item("HOC")
$ws.Activate() | Out-Null

# Core data change: B2 (Opportunity Cost, $/(MW*hr)) goes from 70 to 12.
$ws.Range("B2").Value = 12

# The author's last on-sheet selection moved from F9 to F22 before saving.
$ws.Range("F22").Select() | Out-Null
